# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" sheet (cloned from "2021-Q4" so it keeps the
#   same column layout/formatting) positioned right before "总计".
# - Fills the new sheet with the 2022-Q1 per-fund holding rows.
# - Inserts a new top data row in "总计" summarising 2022-Q1 and renumbers
#   the running index column.

function Set-TextCell($range, $value) {
    # Force the cell to stay text (t="inlineStr"/shared-string "s") even
    # when $value looks numeric (e.g. fund codes like "010695" or ratios
    # like "18.02"), then drop back to the default "Normal" style so we
    # don't leave a stray number-format style behind.
    $range.NumberFormat = "@"
    $range.Value() = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet by cloning "2021-Q4" (same headers /
#    styles), positioned between "2021-Q4" and "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

$fundRows = @(
    @("010695", "华夏磐益一年定期开放混合", "18.02", "82.41", "3.71", "0.6685", 2),
    @("012093", "鹏华创新升级混合型证券投资基金A", "6.58", "64.31", "1.65", "0.1086", 10),
    @("012094", "鹏华创新升级混合型证券投资基金C", "0.19", "64.31", "1.65", "0.0031", 10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    Set-TextCell $q1.Range("B$r") $row[0]
    Set-TextCell $q1.Range("C$r") $row[1]
    Set-TextCell $q1.Range("D$r") $row[2]
    Set-TextCell $q1.Range("E$r") $row[3]
    Set-TextCell $q1.Range("F$r") $row[4]
    Set-TextCell $q1.Range("G$r") $row[5]
    $q1.Range("H$r").Value() = $row[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new row 2 for 2022-Q1 and push the
#    existing quarters down, renumbering the A-column running index.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# restyle/repopulate the freshly inserted row 2
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value() = 0
$total.Range("B2").Value() = "2022-Q1"
$total.Range("C2").Value() = 3
$total.Range("D2").Value() = 0.78

# renumber the running index in column A for the shifted rows
$total.Range("A3").Value() = 1
$total.Range("A4").Value() = 2
$total.Range("A5").Value() = 3

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/selection (unaffected by the diff).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
